$wb = $excel.ActiveWorkbook

# --- Worksheets (by fixed position, independent of name) ---
$wsChains        = $wb.Worksheets.Item(1)   # Basic_chains
$wsCategories    = $wb.Worksheets.Item(2)   # Basic_categories
$wsSubcategories = $wb.Worksheets.Item(3)   # Basic_subcategories
$wsActionTypes   = $wb.Worksheets.Item(4)   # Basic_action_types -> Basic_chain_product_types

# 1) Rename the 4th sheet.
$wsActionTypes.Name = "Basic_chain_product_types"

# 2) Basic_categories!B10 : 10000 -> 1000
$wsCategories.Range("B10").Value = 1000

# 3) Basic_chain_product_types (sheet 4) gets a new empty, formatted cell B4.
#    Copy the format of an existing "empty styled" cell and paste formats only.
$wsChains.Range("C2").Copy()
$wsActionTypes.Range("B4").PasteSpecial(-4122)

# 4) Update selections / active tab.
#    First move the (soon to be non-active) sheet's selection to F22 - this
#    does not change which tab is active yet.
$wsActionTypes.Range("F22").Select()

#    Then activate Basic_categories and select B10, making it the active tab
#    and leaving Basic_chain_product_types as a plain (non-selected) tab.
$wsCategories.Activate()
$wsCategories.Range("B10").Select()
